$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 980
$ws.Range("J7").Value = 980
$ws.Range("L7").Value = 980
$ws.Range("N7").Value = -1204
$ws.Range("H9").Value = 396.92856
$ws.Range("I9").Value = 429.75
$ws.Range("K9").Value = 429.75
$ws.Range("M9").Value = -260.75
$ws.Range("H14").Value = 980
$ws.Range("J14").Value = 980
$ws.Range("L14").Value = 980
$ws.Range("N14").Value = -1362
$ws.Range("H32").Value = 21428.572
$ws.Range("I32").Value = 19000
$ws.Range("J32").Value = 24666.666
$ws.Range("K32").Value = 19000
$ws.Range("L32").Value = 24666.666
$ws.Range("M32").Value = -18674
$ws.Range("N32").Value = -25318.666
$ws.Range("H55").Value = 171.25
$ws.Range("J55").Value = 141.66667
$ws.Range("L55").Value = 141.66667
$ws.Range("N55").Value = -569.6666700000001
$ws.Range("H57").Value = 80000
$ws.Range("J57").Value = 80000
$ws.Range("L57").Value = 240000
$ws.Range("N57").Value = -240998
$ws.Range("H76").Value = 5399
$ws.Range("I76").Value = 5399
$ws.Range("K76").Value = 5399
$ws.Range("M76").Value = -5084
$ws.Range("H79").Value = 5399
$ws.Range("I79").Value = 5399
$ws.Range("K79").Value = 5399
$ws.Range("M79").Value = -4307
$ws.Range("H107").Value = 1173.75
$ws.Range("I107").Value = 1173.75
$ws.Range("K107").Value = 1173.75
$ws.Range("M107").Value = 746.25
$ws.Range("H111").Value = 3966.5
$ws.Range("I111").Value = 4237.25
$ws.Range("J111").Value = 3425
$ws.Range("K111").Value = 12711.75
$ws.Range("L111").Value = 10275
$ws.Range("M111").Value = -9644.75
$ws.Range("N111").Value = -16409
$ws.Range("H127").Value = 2829.6667
$ws.Range("I127").Value = 1795.6
$ws.Range("K127").Value = 5386.799999999999
$ws.Range("M127").Value = -426.7999999999993
$ws.Range("H131").Value = 5796.5713
$ws.Range("I131").Value = 2645
$ws.Range("J131").Value = 9998.666999999999
$ws.Range("K131").Value = 7935
$ws.Range("L131").Value = 29996.001
$ws.Range("M131").Value = -2895
$ws.Range("N131").Value = -40076.001

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3999
$ws.Range("I105").Value = 3999
$ws.Range("K105").Value = 3999
$ws.Range("M105").Value = -2252

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 213.5
$ws.Range("H14").Value = 6998.3335
$ws.Range("I14").Value = 5000
$ws.Range("J14").Value = 7997.5
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 7997.5
$ws.Range("M14").Value = -4830
$ws.Range("N14").Value = -8337.5
$ws.Range("H31").Value = 6356.0386
$ws.Range("I31").Value = 8698.200000000001
$ws.Range("J31").Value = 4892.1875
$ws.Range("K31").Value = 8698.200000000001
$ws.Range("L31").Value = 4892.1875
$ws.Range("M31").Value = -8403.200000000001
$ws.Range("N31").Value = -5482.1875
$ws.Range("H34").Value = 6356.0386
$ws.Range("I34").Value = 8698.200000000001
$ws.Range("J34").Value = 4892.1875
$ws.Range("K34").Value = 8698.200000000001
$ws.Range("L34").Value = 4892.1875
$ws.Range("M34").Value = -8496.200000000001
$ws.Range("N34").Value = -5296.1875
$ws.Range("H62").Value = 3373.75
$ws.Range("I62").Value = 3500
$ws.Range("J62").Value = 2995
$ws.Range("K62").Value = 3500
$ws.Range("L62").Value = 2995
$ws.Range("M62").Value = -2876
$ws.Range("N62").Value = -4243
$ws.Range("H65").Value = 3373.75
$ws.Range("I65").Value = 3500
$ws.Range("J65").Value = 2995
$ws.Range("K65").Value = 17500
$ws.Range("L65").Value = 14975
$ws.Range("M65").Value = -14380
$ws.Range("N65").Value = -21215
$ws.Range("H99").Value = 2546.2
$ws.Range("I99").Value = 2677.3333
$ws.Range("J99").Value = 2349.5
$ws.Range("K99").Value = 2677.3333
$ws.Range("L99").Value = 2349.5
$ws.Range("M99").Value = -1179.3333
$ws.Range("N99").Value = -5345.5
$ws.Range("H126").Value = 2546.2
$ws.Range("I126").Value = 2677.3333
$ws.Range("J126").Value = 2349.5
$ws.Range("K126").Value = 8031.999899999999
$ws.Range("L126").Value = 7048.5
$ws.Range("M126").Value = -5561.999899999999
$ws.Range("N126").Value = -11988.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 12000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = $null
$ws.Range("H108").Value = 2184.3333
$ws.Range("I108").Value = 2184.3333
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 6552.999899999999
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = $null
$ws.Range("N108").Value = -3672.999899999999
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").Value = $null
$ws.Range("H113").Value = 1609.875
$ws.Range("I113").Value = 1654.1428
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 4962.428400000001
$ws.Range("L113").Value = 3900
$ws.Range("M113").Value = -2792.428400000001
$ws.Range("N113").Value = -8240
$ws.Range("H121").Value = 4718
$ws.Range("I121").Value = 850
$ws.Range("J121").Value = 7296.6665
$ws.Range("K121").Value = 2550
$ws.Range("L121").Value = 21889.9995
$ws.Range("M121").Value = -1240
$ws.Range("N121").Value = -24509.9995
$ws.Range("H132").Value = 875.8570999999999
$ws.Range("J132").Value = 846
$ws.Range("L132").Value = 7614
$ws.Range("N132").Value = -12674

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 25000000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null
$ws.Range("H10").Value = 261247.25
$ws.Range("I10").Value = 1000000
$ws.Range("J10").Value = 14996.333
$ws.Range("K10").Value = 1000000
$ws.Range("L10").Value = 14996.333
$ws.Range("M10").Value = -999831
$ws.Range("N10").Value = -15334.333
$ws.Range("H15").Value = 26666.666
$ws.Range("I15").Value = 20000
$ws.Range("K15").Value = 20000
$ws.Range("M15").Value = -19712
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = $null
$ws.Range("N34").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = $null
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = $null
$ws.Range("N79").Value = 0
$ws.Range("H81").Value = 26666.666
$ws.Range("I81").Value = 20000
$ws.Range("K81").Value = 20000
$ws.Range("M81").Value = -19002
$ws.Range("H84").Value = 26666.666
$ws.Range("I84").Value = 20000
$ws.Range("K84").Value = 60000
$ws.Range("M84").Value = -55008
$ws.Range("H109").Value = 47500
$ws.Range("J109").Value = 47500
$ws.Range("L109").Value = 47500
$ws.Range("N109").Value = -49580
$ws.Range("H132").Value = 15153638
$ws.Range("I132").Value = 1570.0588
$ws.Range("K132").Value = 4710.1764
$ws.Range("M132").Value = -2180.1764

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 41667380
$ws.Range("I22").Value = 492.6
$ws.Range("J22").Value = 71429440
$ws.Range("K22").Value = 492.6
$ws.Range("L22").Value = 71429440
$ws.Range("M22").Value = -197.6
$ws.Range("N22").Value = -71430030
$ws.Range("H27").Value = 41667380
$ws.Range("I27").Value = 492.6
$ws.Range("J27").Value = 71429440
$ws.Range("K27").Value = 492.6
$ws.Range("L27").Value = 71429440
$ws.Range("M27").Value = -385.6
$ws.Range("N27").Value = -71429654
$ws.Range("H40").Value = 83338950
$ws.Range("I40").Value = 142861360
$ws.Range("K40").Value = 142861360
$ws.Range("M40").Value = -142861224
$ws.Range("H46").Value = 3774.2666
$ws.Range("I46").Value = 1719.8334
$ws.Range("K46").Value = 1719.8334
$ws.Range("M46").Value = -1531.8334
$ws.Range("H68").Value = 1018.8333
$ws.Range("I68").Value = 1062.6
$ws.Range("J68").Value = 800
$ws.Range("K68").Value = 1062.6
$ws.Range("L68").Value = 800
$ws.Range("M68").Value = -313.5999999999999
$ws.Range("N68").Value = -2298
$ws.Range("H71").Value = 1018.8333
$ws.Range("I71").Value = 1062.6
$ws.Range("J71").Value = 800
$ws.Range("K71").Value = 5313
$ws.Range("L71").Value = 4000
$ws.Range("M71").Value = -1569
$ws.Range("N71").Value = -11488
$ws.Range("H128").Value = 73153.57000000001
$ws.Range("J128").Value = 73153.57000000001
$ws.Range("L128").Value = 73153.57000000001
$ws.Range("N128").Value = -83113.57000000001
$ws.Range("H136").Value = 250000850
$ws.Range("I136").Value = 1698
$ws.Range("K136").Value = 5094
$ws.Range("M136").Value = -2544

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").Value = $null
